$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "TC ID"
$ws.Range("B1").Value = "TCD ID"
$ws.Range("C1").Value = "TR ID"
